# Generate Report for Handoff
#
# The localization-status report is regenerated: the handoff run for the
# "06d00e15-5718-4b72-ba2e-f3c31e04ba4a" file (and the other files that
# share its handoff batch / handback timestamp) now carries the "ht"
# (handoff type) priority, and the handoff/handback timestamps for that
# batch move forward a few seconds to reflect the new run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" (Overview!G) / "Latest Handback DateTime"
# (de-de!H) both held the same shared timestamp string; replace it on both
# sheets so every cell that displayed it moves to the new value.
[void]$wsOverview.Cells.Replace("2016-08-26 00:21:22", "2016-08-26 00:21:39")
[void]$wsDeDe.Cells.Replace("2016-08-26 00:21:22", "2016-08-26 00:21:39")

# zh-cn!H held its own shared timestamp string ("Latest Handback DateTime").
[void]$wsZhCn.Cells.Replace("2016-08-26 00:21:18", "2016-08-26 00:21:34")

# The rows for 06d00e15, 12a67c01, 37935410, 5ec2ccaf, c2e9f80d and
# cfabdd89 (rows 7, 9, 10, 12, 13, 14) now report a "Priority" of "ht" on
# both the zh-cn and de-de handback sheets.
foreach ($r in 7, 9, 10, 12, 13, 14) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}
